# Applies the crypto-price/volume refresh described by the commit:
#   "Updated cryptos list on Sun Feb 25 17:37:07 UTC 2024 with GitHub Actions"
#
# Price (column D) and Volume(1h) (column E) text is refreshed for most rows,
# and rows 42/43 (Stellar / Stacks) swap places entirely (name, link, price,
# volume all move up/down one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates -----------------------------------------------------
# These source values are unambiguous as text (they contain a percent sign,
# surrounding spaces, letters, or more than one "." separator) so a normal
# Range.Value assignment keeps them stored as text, matching the workbook's
# original inline-string cells without touching any cell style.
$ws.Range('D2').Value = '51.471.29'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '3.051.99'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '3.529.90'
$ws.Range('E13').Value = '  +2.22%  '
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').Value = '3.043.48'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('E18').Value = '  -5.70%  '
$ws.Range('D19').Value = '51.529.18'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('E26').Value = '  +5.35%  '
$ws.Range('E27').Value = '  +3.09%  '
$ws.Range('E28').Value = '  +3.00%  '
$ws.Range('E29').Value = '  -3.88%  '
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +2.54%  '
$ws.Range('E39').Value = '  +8.02%  '
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('E41').Value = '  +1.25%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('E47').Value = '  +3.52%  '
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('D49').Value = '2.028.16'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').Value = '3.350.13'
$ws.Range('E50').Value = '  +2.37%  '
$ws.Range('E51').Value = '  -3.81%  '

# --- Numeric-looking text updates --------------------------------------------
# These new prices (e.g. "0.999", "2.40") are valid numeric literals, so a
# direct Range.Value assignment would make Excel auto-convert them to the
# Number type (and could introduce binary floating-point noise). The source
# column is plain untyped text, so instead we render the exact text via a
# TEXT() formula in a scratch cell and bring over only the computed value with
# PasteSpecial (xlPasteValues = -4163), which keeps the destination cell a
# literal string and leaves its style untouched.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '=TEXT(0.999,"0.000")'
$scratch.Copy()
$ws.Range('D4').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(385.09,"0.00")'
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(103.06,"0.00")'
$scratch.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(0.585,"0.000")'
$scratch.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(36.72,"0.00")'
$scratch.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(0.0862,"0.0000")'
$scratch.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(18.53,"0.00")'
$scratch.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(0.971,"0.000")'
$scratch.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(10.63,"0.00")'
$scratch.Copy()
$ws.Range('D18').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(12.41,"0.00")'
$scratch.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(70.11,"0.00")'
$scratch.Copy()
$ws.Range('D23').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(267.77,"0.00")'
$scratch.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(8.23,"0.00")'
$scratch.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(26.88,"0.00")'
$scratch.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(7.21,"0.00")'
$scratch.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(10.27,"0.00")'
$scratch.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(34.61,"0.00")'
$scratch.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(50.35,"0.00")'
$scratch.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(0.0447,"0.0000")'
$scratch.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(2.55,"0.00")'
$scratch.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(0.116,"0.000")'
$scratch.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(124.92,"0.00")'
$scratch.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(21.85,"0.00")'
$scratch.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$scratch.Formula = '=TEXT(2.40,"0.00")'
$scratch.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$scratch.Value = ""
$excel.CutCopyMode = 0

